# Sync attendance_reports: swap the order of the two comma-separated
# "Recorded By" names in column G so that "dnasr281@gmail.com, X"
# becomes "X, dnasr281@gmail.com". Cells containing only a single
# name (no comma) or three names are left untouched, matching the
# source diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $text = $cell.Value2

    if ($text -ne $null -and $text -like "dnasr281@gmail.com, *") {
        $rest = $text.Substring("dnasr281@gmail.com, ".Length)
        if ($rest -notlike "*,*") {
            $newText = $rest + ", dnasr281@gmail.com"
            $cell.Value = $newText
        }
    }
}
